$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.764.48'
$ws.Range('E2').Value = '  -0.21%  '
$ws.Range('D3').Value = '1.595.21'
$ws.Range('E3').Value = '  -1.48%  '
$ws.Range('E4').Value = '  +0.22%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '209.79'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.59%  '
$ws.Range('E6').Value = '  -1.92%  '
$ws.Range('E7').Value = '  +0.20%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '22.37'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.50%  '
$ws.Range('E9').Value = '  -1.23%  '
$ws.Range('E10').Value = '  -1.63%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0870'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.26%  '
$ws.Range('D12').Value = '1.822.19'
$ws.Range('E12').Value = '  -1.53%  '
$ws.Range('D13').Value = '1.606.19'
$ws.Range('E13').Value = '  -1.77%  '
$ws.Range('E14').Value = '  -2.36%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.533'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.37%  '
$ws.Range('D16').Value = '27.765.37'
$ws.Range('E16').Value = '  -0.15%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.49'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.50%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '219.36'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.90%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.39'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.49%  '
$ws.Range('E20').Value = '  -2.20%  '
$ws.Range('E21').Value = '  +0.27%  '
$ws.Range('E22').Value = '  -3.17%  '
$ws.Range('E23').Value = '  -1.94%  '
$ws.Range('E24').Value = '  -4.00%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '154.08'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.31%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.18'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +4.31%  '
$ws.Range('E27').Value = '  +0.21%  '
$ws.Range('E28').Value = '  -0.78%  '
$ws.Range('E29').Value = '  -3.39%  '
$ws.Range('E30').Value = '  -0.72%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0474'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.97%  '
$ws.Range('E32').Value = '  -3.91%  '
$ws.Range('D33').Value = '1.379.70'
$ws.Range('E33').Value = '  -1.90%  '
$ws.Range('E34').Value = '  -2.57%  '
$ws.Range('E35').Value = '  -3.57%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.977'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.79%  '
$ws.Range('E37').Value = '  +0.15%  '
$ws.Range('E39').Value = '  -2.44%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.829'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.55%  '
$ws.Range('E41').Value = '  +0.16%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.975'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.10%  '
$ws.Range('E43').Value = '  -1.01%  '
$ws.Range('E44').Value = '  +2.73%  '
$ws.Range('B45').Value = 'RenderToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.75'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.61%  '
$ws.Range('B46').Value = 'FraxShare'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '5.23'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.09%  '
$ws.Range('D47').Value = '1.733.03'
$ws.Range('E47').Value = '  -1.49%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '86.47'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.34%  '
$ws.Range('E49').Value = '  -0.92%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0967'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.47%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0496'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.05%  '
